$d = $word.ActiveDocument

# Locate the "Ver no Jupiter ..." paragraph; this anchors the block of
# trailing site-chrome paragraphs (a blank paragraph, the "Ver no
# Jupiter ..." line, and the "© 2020 ..." line) that the site rebuild
# dropped from the bibliography section.
$anchor = $d.Content
$found = $anchor.Find.Execute("Ver no Jupiter Salvar em pdf Salvar em docx", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if (-not $found) {
    throw "Anchor paragraph 'Ver no Jupiter ...' not found"
}

$anchorIndex = $anchor.Paragraphs.First.Index

# Paragraph right before the anchor is the blank separator paragraph;
# paragraph right after it is the "© 2020 ..." copyright line. Delete
# the whole span (blank + "Ver no Jupiter ..." + "© 2020 ...") in one
# go, leaving the bibliography paragraph and the remaining blank /
# page-break paragraphs untouched.
$startPara = $d.Paragraphs.Item($anchorIndex - 1)
$endPara = $d.Paragraphs.Item($anchorIndex + 1)

$deleteRange = $d.Range($startPara.Range.Start, $endPara.Range.End)
$deleteRange.Delete()
